# Apply the diff: swap row 28 ("new" kitchens item) and row 29 ("catch" item)
# such that the "new" entry (with its M:V stats) moves down to row 29, and
# row 28 becomes a "catch" entry pointing at a different stimulus image.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: turn the old "catch" row into the "new" kitchens row ---
# Add H29 (category = kitchens), which row 29 previously lacked.
$ws.Range("H29").Value = "kitchens"

# J29: catch -> new
$ws.Range("J29").Value = "new"

# L29: stimuli/catch_09_stairs.jpg -> stimuli/img_57os5.png
$ws.Range("L29").Value = "stimuli/img_57os5.png"

# Add the statistic columns M29:V29 that row 28 used to carry.
$ws.Range("M29").Value = 82.70588235294117
$ws.Range("N29").Value = 65.73529411764706
$ws.Range("O29").Value = 74.22058823529412
$ws.Range("P29").Value = 34
$ws.Range("Q29").Value = 9
$ws.Range("R29").Value = 9
$ws.Range("S29").Value = 9
$ws.Range("T29").Value = 9
$ws.Range("U29").Value = 9
$ws.Range("V29").Value = 9

# --- Row 28: turn the old "new" kitchens row into the "catch" row ---
# Remove H28 (category), which the new "catch" row does not have.
$ws.Range("H28").ClearContents()

# J28: new -> catch
$ws.Range("J28").Value = "catch"

# L28: stimuli/img_57os5.png -> stimuli/catch_16.jpg
$ws.Range("L28").Value = "stimuli/catch_16.jpg"

# Remove the statistic columns M28:V28 since "catch" rows don't carry them.
$ws.Range("M28:V28").ClearContents()
